# Refresh the cryptos price table with the latest scrape values.
# (GitHub Actions scheduled update: prices/volumes change, and the
#  WrappedEther / Polkadot rows (12 & 13) swap places.)
#
# NumberFormat "@" (Text) is applied before each write, and the style is
# reset back to "Normal" afterwards, so purely-numeric-looking strings
# (e.g. "9.960", "0.5540") are stored verbatim as text -- matching the
# original inline-string cells -- instead of being coerced into numbers
# and losing significant trailing zeros / formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "25.888.54"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -0.57%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.640.76"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.49%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "214.88"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.19%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5049"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +0.01%  "
$cell.Style = "Normal"

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.53%  "
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2567"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -0.20%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06387"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -0.65%  "
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.58"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +0.74%  "
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07796"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +1.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("B12")
$cell.NumberFormat = "@"
$cell.Value = "WrappedEther"
$cell.Style = "Normal"
$cell = $ws.Range("C12")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.654.90"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +0.89%  "
$cell.Style = "Normal"

$cell = $ws.Range("B13")
$cell.NumberFormat = "@"
$cell.Value = "Polkadot"
$cell.Style = "Normal"
$cell = $ws.Range("C13")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.274"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +0.79%  "
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.5422"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -0.39%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅7858"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -0.48%  "
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "64.68"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +1.75%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "25.935.61"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -0.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -0.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "198.08"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -2.79%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "4.385"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +2.44%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "9.960"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -0.19%  "
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.983"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +0.38%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -0.48%  "
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.862"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -3.40%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "140.14"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -0.81%  "
$cell.Style = "Normal"

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.1142"
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "6.844"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  +1.71%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "15.72"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +0.11%  "
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.242"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.04934"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -2.26%  "
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.262"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +0.51%  "
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.193"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.26%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.531"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -0.56%  "
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "2.366"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +1.11%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.8924"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.604"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -0.61%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.139.24"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -0.43%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5540"
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01558"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -0.81%  "
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -0.46%  "
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.668"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.8200"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +1.31%  "
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "99.23"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -0.45%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0₈119"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +5.05%  "
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.778.41"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.4517"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -0.28%  "
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "55.25"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +0.82%  "
$cell.Style = "Normal"

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -0.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.05051"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +0.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.006"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -0.22%  "
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.09507"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +2.06%  "
$cell.Style = "Normal"
